$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.593.03'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.600.52'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.39'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.83'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.49%  '
$ws.Range("E9").Value = '  +1.18%  '
$ws.Range("E10").Value = '  +1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0912'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.829.70'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.606.79'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.540'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +3.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.607.50'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("E16").Value = '  +0.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.81'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '241.69'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.65'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0693'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.24'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.14'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.36'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.39%  '
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("E30").Value = '  +2.31%  '
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.15'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.422.71'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  +2.47%  '
$ws.Range("E36").Value = '  +4.37%  '
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.30'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("E39").Value = '  +3.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.544'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '55.78'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.88%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("E43").Value = '  +5.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.810'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.78%  '
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("E46").Value = '  +16.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '66.36'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.81%  '
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.740.36'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '86.08'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("E51").Value = '  +3.71%  '
